$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "528.27", "0.999") that must
# stay plain text, exactly like the original inlineStr cells. Excel's COM
# layer auto-coerces such strings assigned via .Value into real numbers,
# so for each D cell we briefly force a text number format, assign the
# value, then clear the format again so no stray cell style is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.001.81"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.342.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.83"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.339.95"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.51%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.609"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.135"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.873.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.337.69"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.84%  "

$ws.Range("E17").Value = "  +1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.924.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.964"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.62"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +9.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.72"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "634.41"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.23"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("E35").Value = "  +1.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.72"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.48"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.381"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0725"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +13.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +10.29%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.970.75"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.13%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.125"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.03"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.70"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0396"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.04"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.125"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.11"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.36%  "
